$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.68314440705249
$ws.Cells.Item(2, 3).Value = 7.44904072262006
$ws.Cells.Item(2, 4).Value = 11.48108678349906
$ws.Cells.Item(2, 6).Value = 28.07982124531163
$ws.Cells.Item(2, 7).Value = 25.31283876543285
$ws.Cells.Item(2, 8).Value = 13.34986025467813
$ws.Cells.Item(2, 10).Value = 11.21989229887501
$ws.Cells.Item(2, 13).Value = 16.49252308586409
$ws.Cells.Item(2, 14).Value = 17.03302699124378
$ws.Cells.Item(2, 15).Value = 19.87165092779136
$ws.Cells.Item(3, 2).Value = 11.17585858336026
$ws.Cells.Item(3, 3).Value = 7.118073954361782
$ws.Cells.Item(3, 4).Value = 11.46439291866277
$ws.Cells.Item(3, 6).Value = 28.10129476322183
$ws.Cells.Item(3, 7).Value = 25.28473549344122
$ws.Cells.Item(3, 8).Value = 13.38965292054986
$ws.Cells.Item(3, 10).Value = 11.24722677181777
$ws.Cells.Item(3, 13).Value = 16.33357170033008
$ws.Cells.Item(3, 14).Value = 17.07569612195835
$ws.Cells.Item(3, 15).Value = 19.92072144201968
$ws.Cells.Item(4, 2).Value = 10.85320893001702
$ws.Cells.Item(4, 3).Value = 6.905726496568719
$ws.Cells.Item(4, 4).Value = 11.45602244036113
$ws.Cells.Item(4, 6).Value = 28.12196722358647
$ws.Cells.Item(4, 7).Value = 25.2770065025454
$ws.Cells.Item(4, 8).Value = 13.41645970570025
$ws.Cells.Item(4, 10).Value = 11.26562955622145
$ws.Cells.Item(4, 13).Value = 16.23763683849317
$ws.Cells.Item(4, 14).Value = 17.10362961990512
$ws.Cells.Item(4, 15).Value = 19.95573983319532
$ws.Cells.Item(5, 2).Value = 10.71909691433692
$ws.Cells.Item(5, 3).Value = 6.816977033231961
$ws.Cells.Item(5, 4).Value = 11.45308691128494
$ws.Cells.Item(5, 6).Value = 28.13227203604722
$ws.Cells.Item(5, 7).Value = 25.27625281819564
$ws.Cells.Item(5, 8).Value = 13.42798011169156
$ws.Cells.Item(5, 10).Value = 11.27353600884418
$ws.Cells.Item(5, 13).Value = 16.19899754528203
$ws.Cells.Item(5, 14).Value = 17.11544963334899
$ws.Cells.Item(5, 15).Value = 19.97123629949267
$ws.Cells.Item(6, 2).Value = 10.69667463693277
$ws.Cells.Item(6, 3).Value = 6.802108994072434
$ws.Cells.Item(6, 4).Value = 11.45262826612207
$ws.Cells.Item(6, 6).Value = 28.13409664270634
$ws.Cells.Item(6, 7).Value = 25.27627233647843
$ws.Cells.Item(6, 8).Value = 13.42992907316408
$ws.Cells.Item(6, 10).Value = 11.27487346045108
$ws.Cells.Item(6, 13).Value = 16.19261006383365
$ws.Cells.Item(6, 14).Value = 17.11743874635916
$ws.Cells.Item(6, 15).Value = 19.97388343660974
$ws.Cells.Item(7, 2).Value = 10.85141063990601
$ws.Cells.Item(7, 3).Value = 6.904538448178071
$ws.Cells.Item(7, 4).Value = 11.45598092191534
$ws.Cells.Item(7, 6).Value = 28.12209858702236
$ws.Cells.Item(7, 7).Value = 25.27698663853675
$ws.Cells.Item(7, 8).Value = 13.41661265974042
$ws.Cells.Item(7, 10).Value = 11.26573453683123
$ws.Cells.Item(7, 13).Value = 16.23711384510671
$ws.Cells.Item(7, 14).Value = 17.10378725881927
$ws.Cells.Item(7, 15).Value = 19.95594386357515
$ws.Cells.Item(8, 2).Value = 11.51065213777613
$ws.Cells.Item(8, 3).Value = 7.336863276778803
$ws.Cells.Item(8, 4).Value = 11.47494242881909
$ws.Cells.Item(8, 6).Value = 28.08567022245957
$ws.Cells.Item(8, 7).Value = 25.30117261203446
$ws.Cells.Item(8, 8).Value = 13.36308775252538
$ws.Cells.Item(8, 10).Value = 11.22898113024684
$ws.Cells.Item(8, 13).Value = 16.4373941628355
$ws.Cells.Item(8, 14).Value = 17.04737975731909
$ws.Cells.Item(8, 15).Value = 19.88755392957974
$ws.Cells.Item(9, 2).Value = 12.70801949261509
$ws.Cells.Item(9, 3).Value = 8.109248572027946
$ws.Cells.Item(9, 4).Value = 11.52690207487337
$ws.Cells.Item(9, 6).Value = 28.0737019701161
$ws.Cells.Item(9, 7).Value = 25.42404301081082
$ws.Cells.Item(9, 8).Value = 13.27698417647919
$ws.Cells.Item(9, 10).Value = 11.16975697901931
$ws.Cells.Item(9, 13).Value = 16.84167106628555
$ws.Cells.Item(9, 14).Value = 16.9504943352803
$ws.Cells.Item(9, 15).Value = 19.79236301856312
$ws.Cells.Item(10, 2).Value = 13.52200370067978
$ws.Cells.Item(10, 3).Value = 8.627530136155249
$ws.Cells.Item(10, 4).Value = 11.57388274243842
$ws.Cells.Item(10, 6).Value = 28.10117716345722
$ws.Cells.Item(10, 7).Value = 25.55992216640944
$ws.Cells.Item(10, 8).Value = 13.22525035311167
$ws.Cells.Item(10, 10).Value = 11.13407881746466
$ws.Cells.Item(10, 13).Value = 17.14346866248732
$ws.Cells.Item(10, 14).Value = 16.88763780066706
$ws.Cells.Item(10, 15).Value = 19.74631508508885
$ws.Cells.Item(11, 2).Value = 13.87680463132691
$ws.Cells.Item(11, 3).Value = 8.852121168660352
$ws.Cells.Item(11, 4).Value = 11.59711833691075
$ws.Cells.Item(11, 6).Value = 28.12153221935067
$ws.Cells.Item(11, 7).Value = 25.6314980544671
$ws.Cells.Item(11, 8).Value = 13.20422401161514
$ws.Cells.Item(11, 10).Value = 11.11954882586654
$ws.Cells.Item(11, 13).Value = 17.28131878176043
$ws.Cells.Item(11, 14).Value = 16.86084127972057
$ws.Cells.Item(11, 15).Value = 19.73058177003163
$ws.Cells.Item(12, 2).Value = 14.00884971397662
$ws.Cells.Item(12, 3).Value = 8.935527456334986
$ws.Cells.Item(12, 4).Value = 11.60618039912634
$ws.Cells.Item(12, 6).Value = 28.13036629009335
$ws.Cells.Item(12, 7).Value = 25.65998961389563
$ws.Cells.Item(12, 8).Value = 13.19662292338122
$ws.Cells.Item(12, 10).Value = 11.11429115446033
$ws.Cells.Item(12, 13).Value = 17.33355938552639
$ws.Cells.Item(12, 14).Value = 16.85095190766522
$ws.Cells.Item(12, 15).Value = 19.72537543571053
$ws.Cells.Item(13, 2).Value = 13.98051538794408
$ws.Cells.Item(13, 3).Value = 8.917637904341069
$ws.Cells.Item(13, 4).Value = 11.60421709296081
$ws.Cells.Item(13, 6).Value = 28.12841370462098
$ws.Cells.Item(13, 7).Value = 25.65379204021142
$ws.Cells.Item(13, 8).Value = 13.19824388420064
$ws.Cells.Item(13, 10).Value = 11.11541261255214
$ws.Cells.Item(13, 13).Value = 17.3223073151863
$ws.Cells.Item(13, 14).Value = 16.8530702997885
$ws.Cells.Item(13, 15).Value = 19.72646326465641
$ws.Cells.Item(14, 2).Value = 13.88771479455292
$ws.Cells.Item(14, 3).Value = 8.859016125641487
$ws.Cells.Item(14, 4).Value = 11.59785863072222
$ws.Cells.Item(14, 6).Value = 28.12223646119593
$ws.Cells.Item(14, 7).Value = 25.63381437356434
$ws.Cells.Item(14, 8).Value = 13.2035914249724
$ws.Cells.Item(14, 10).Value = 11.11911137336765
$ws.Cells.Item(14, 13).Value = 17.2856160836324
$ws.Cells.Item(14, 14).Value = 16.86002250974554
$ws.Cells.Item(14, 15).Value = 19.73013837132293
$ws.Cells.Item(15, 2).Value = 13.83056862154302
$ws.Cells.Item(15, 3).Value = 8.822893916876943
$ws.Cells.Item(15, 4).Value = 11.5939980249571
$ws.Cells.Item(15, 6).Value = 28.11859923700891
$ws.Cells.Item(15, 7).Value = 25.62175760397047
$ws.Cells.Item(15, 8).Value = 13.20691399328343
$ws.Cells.Item(15, 10).Value = 11.12140881578808
$ws.Cells.Item(15, 13).Value = 17.26314561565979
$ws.Cells.Item(15, 14).Value = 16.86431450777958
$ws.Cells.Item(15, 15).Value = 19.73248739712299
$ws.Cells.Item(16, 2).Value = 13.49849757666911
$ws.Cells.Item(16, 3).Value = 8.61262460286
$ws.Cells.Item(16, 4).Value = 11.57240134578589
$ws.Cells.Item(16, 6).Value = 28.10000471527812
$ws.Cells.Item(16, 7).Value = 25.55543962196627
$ws.Cells.Item(16, 8).Value = 13.22667498067972
$ws.Cells.Item(16, 10).Value = 11.1350626022063
$ws.Cells.Item(16, 13).Value = 17.13446761408521
$ws.Cells.Item(16, 14).Value = 16.88942513373525
$ws.Cells.Item(16, 15).Value = 19.74744837355448
$ws.Cells.Item(17, 2).Value = 13.29075291267364
$ws.Cells.Item(17, 3).Value = 8.480742046869644
$ws.Cells.Item(17, 4).Value = 11.55962655643611
$ws.Cells.Item(17, 6).Value = 28.09060754466768
$ws.Cells.Item(17, 7).Value = 25.51724542876117
$ws.Cells.Item(17, 8).Value = 13.23944036499572
$ws.Cells.Item(17, 10).Value = 11.1438742479061
$ws.Cells.Item(17, 13).Value = 17.05564226356156
$ws.Cells.Item(17, 14).Value = 16.90528960605124
$ws.Cells.Item(17, 15).Value = 19.75796327112033
$ws.Cells.Item(18, 2).Value = 13.16981027655523
$ws.Cells.Item(18, 3).Value = 8.403836195527914
$ws.Cells.Item(18, 4).Value = 11.55245457446394
$ws.Cells.Item(18, 6).Value = 28.08594232212366
$ws.Cells.Item(18, 7).Value = 25.49619733009175
$ws.Cells.Item(18, 8).Value = 13.24701867703818
$ws.Cells.Item(18, 10).Value = 11.14910250172335
$ws.Cells.Item(18, 13).Value = 17.0103592274399
$ws.Cells.Item(18, 14).Value = 16.91458362669532
$ws.Cells.Item(18, 15).Value = 19.76450185883912
$ws.Cells.Item(19, 2).Value = 13.12861432669071
$ws.Cells.Item(19, 3).Value = 8.377617912081194
$ws.Cells.Item(19, 4).Value = 11.55005658861673
$ws.Cells.Item(19, 6).Value = 28.08448990965971
$ws.Cells.Item(19, 7).Value = 25.48922931255602
$ws.Cells.Item(19, 8).Value = 13.24962507967862
$ws.Cells.Item(19, 10).Value = 11.15090018220668
$ws.Cells.Item(19, 13).Value = 16.99503790620443
$ws.Cells.Item(19, 14).Value = 16.91775949901163
$ws.Cells.Item(19, 15).Value = 19.7667999305846
$ws.Cells.Item(20, 2).Value = 13.31301875671402
$ws.Cells.Item(20, 3).Value = 8.494890158385768
$ws.Cells.Item(20, 4).Value = 11.56096829990562
$ws.Cells.Item(20, 6).Value = 28.09153134797947
$ws.Cells.Item(20, 7).Value = 25.52121614165965
$ws.Cells.Item(20, 8).Value = 13.23805703885448
$ws.Cells.Item(20, 10).Value = 11.14291967052421
$ws.Cells.Item(20, 13).Value = 17.06402793910415
$ws.Cells.Item(20, 14).Value = 16.90358329880521
$ws.Cells.Item(20, 15).Value = 19.75679314319422
$ws.Cells.Item(21, 2).Value = 13.9150358814268
$ws.Cells.Item(21, 3).Value = 8.876279554093918
$ws.Cells.Item(21, 4).Value = 11.59971916022961
$ws.Cells.Item(21, 6).Value = 28.12402034179946
$ws.Cells.Item(21, 7).Value = 25.63964479362427
$ws.Cells.Item(21, 8).Value = 13.20201091754461
$ws.Cells.Item(21, 10).Value = 11.11801832122881
$ws.Cells.Item(21, 13).Value = 17.29639243258316
$ws.Cells.Item(21, 14).Value = 16.85797348243006
$ws.Cells.Item(21, 15).Value = 19.72903849504217
$ws.Cells.Item(22, 2).Value = 14.2949926599543
$ws.Cells.Item(22, 3).Value = 9.115960093638149
$ws.Cells.Item(22, 4).Value = 11.62657707759676
$ws.Cells.Item(22, 6).Value = 28.15181473926275
$ws.Cells.Item(22, 7).Value = 25.72512148773701
$ws.Cells.Item(22, 8).Value = 13.18055786283055
$ws.Cells.Item(22, 10).Value = 11.10316910900589
$ws.Cells.Item(22, 13).Value = 17.44846915993382
$ws.Cells.Item(22, 14).Value = 16.82966776133802
$ws.Cells.Item(22, 15).Value = 19.71528026428442
$ws.Cells.Item(23, 2).Value = 14.09346178778785
$ws.Cells.Item(23, 3).Value = 8.988924364082322
$ws.Cells.Item(23, 4).Value = 11.61210397484945
$ws.Cells.Item(23, 6).Value = 28.13638153939721
$ws.Cells.Item(23, 7).Value = 25.67876816728201
$ws.Cells.Item(23, 8).Value = 13.19181497144702
$ws.Cells.Item(23, 10).Value = 11.11096399465607
$ws.Cells.Item(23, 13).Value = 17.36729684058247
$ws.Cells.Item(23, 14).Value = 16.84463771660826
$ws.Cells.Item(23, 15).Value = 19.72222194724966
$ws.Cells.Item(24, 2).Value = 13.30295705753665
$ws.Cells.Item(24, 3).Value = 8.488497173502342
$ws.Cells.Item(24, 4).Value = 11.56036116027624
$ws.Cells.Item(24, 6).Value = 28.09111139948485
$ws.Cells.Item(24, 7).Value = 25.51941814575692
$ws.Cells.Item(24, 8).Value = 13.23868169516102
$ws.Cells.Item(24, 10).Value = 11.14335072934347
$ws.Cells.Item(24, 13).Value = 17.06023666395024
$ws.Cells.Item(24, 14).Value = 16.90435418036963
$ws.Cells.Item(24, 15).Value = 19.75732062114713
$ws.Cells.Item(25, 2).Value = 12.39516867146766
$ws.Cells.Item(25, 3).Value = 7.908738370784415
$ws.Cells.Item(25, 4).Value = 11.511283353835
$ws.Cells.Item(25, 6).Value = 28.0705673488296
$ws.Cells.Item(25, 7).Value = 25.38275514437422
$ws.Cells.Item(25, 8).Value = 13.29825533992685
$ws.Cells.Item(25, 10).Value = 11.18440296247369
$ws.Cells.Item(25, 13).Value = 16.73130081276203
$ws.Cells.Item(25, 14).Value = 16.97523925797426
$ws.Cells.Item(25, 15).Value = 19.81392992510473
